# Applies the "nuevos casos de uso" edit described in the commit message:
#  1. Cesion de contrato nit a nit con cambio a plan pospago empresarial 5.3
#  2. Activacion nintendo con cliente nit
#
# Changes touch the "Semilla 9" sheet (selection only) and the "Semilla 8"
# sheet (the actual new data rows / selection / page setup).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Semilla 9")
$ws3 = $wb.Worksheets.Item("Semilla 8")

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to stay a text value even when the string looks like
    # a number (e.g. "3045981670"), matching how these sheets already
    # store MSISDN / IMEI-like codes as shared strings rather than numbers.
    $range.NumberFormat = "@"
    $range.Value2 = $text
}

# ---------------------------------------------------------------------
# Sheet "Semilla 9": only the active selection changed in the source file.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F2:G2").Select()

# ---------------------------------------------------------------------
# Sheet "Semilla 8": new use-case rows.
# ---------------------------------------------------------------------

# Row 9 - caso "cesion de contrato nit a nit" (update)
Set-TextValue $ws3.Range("B9") "459399130"
Set-TextValue $ws3.Range("C9") "3052749177"
Set-TextValue $ws3.Range("D9") "732111193280551"

# Row 10 (update)
Set-TextValue $ws3.Range("B10") "836898669"
Set-TextValue $ws3.Range("C10") "3052754285"
Set-TextValue $ws3.Range("D10") "732111324709512"

# Row 11 (update)
Set-TextValue $ws3.Range("B11") "255188531"
Set-TextValue $ws3.Range("C11") "3052749177"
Set-TextValue $ws3.Range("D11") "732111193280551"

# Row 12 (update) + new column E
Set-TextValue $ws3.Range("B12") "194936717"
Set-TextValue $ws3.Range("C12") "3045987650"
Set-TextValue $ws3.Range("D12") "732111324709673"
Set-TextValue $ws3.Range("E12") "client nit a nit"

# Row 13 (update) + new column E
Set-TextValue $ws3.Range("B13") "432694001"
Set-TextValue $ws3.Range("C13") "3046010569"
Set-TextValue $ws3.Range("D13") "732111324709674"
Set-TextValue $ws3.Range("E13") "908348697"

# Row 14: B14 used to hold a plain number (920626579); it becomes a text
# value left-aligned, same as column C already was.
$ws3.Range("B14").HorizontalAlignment = -4131
Set-TextValue $ws3.Range("B14") "920626579"
Set-TextValue $ws3.Range("C14") "3045984556"
Set-TextValue $ws3.Range("D14") "732111324709675"

# New row 15 - "activacion nintendo con cliente nit"
Set-TextValue $ws3.Range("A15") "10960370"
Set-TextValue $ws3.Range("B15") "211423738"
Set-TextValue $ws3.Range("C15") "3052754289"
Set-TextValue $ws3.Range("D15") "732111324709676"

# New row 16
Set-TextValue $ws3.Range("A16") "10960370"
Set-TextValue $ws3.Range("B16") "971449311"
Set-TextValue $ws3.Range("C16") "3046008586"
Set-TextValue $ws3.Range("D16") "732111193278871"

# Page is now printed in portrait orientation.
$ws3.PageSetup.Orientation = 1

# Active selection moves past the new rows.
$ws3.Activate()
$ws3.Range("A17").Select()
